$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (row 8), columns A-F first ---
$ws.Range("A8").Value = "PaymentEmail"
$ws.Range("B8").Value = "PaymentCardNo"
$ws.Range("C8").Value = "PaymentMMYY"
$ws.Range("D8").Value = "Paymentcvv"
$ws.Range("E8").Value = "PaymentNameOnCard"
$ws.Range("F8").Value = "PaymentCountry"

# --- Data row (row 9): name/country before the remaining headers, to
# reproduce the original shared-string insertion order ---
$ws.Range("E9").Value = "Mahesh"
$ws.Range("F9").Value = "India"

# --- Remaining header cells (G8:H8) ---
$ws.Range("G8").Value = "Extended_time"
$ws.Range("H8").Value = "Extended date"

# --- Data row (row 9): date / numeric cells ---
$ws.Range("C9").NumberFormat = "mm-dd-yy"
$ws.Range("C9").Value = 44585

$ws.Range("B9").NumberFormat = "0"
$ws.Range("B9").Value = 4242424242424240

$ws.Range("D9").Value = 234

# Copy C9's date format onto G9 / H9 so they share the exact same style
$ws.Range("C9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("H9").PasteSpecial(-4122)
$ws.Range("H9").Value = 44677

# --- Hyperlinked consumer e-mail (A9) ---
$ws.Range("A9").Value = "consumer2@gmail.com"
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:consumer2@gmail.com", "", "", "consumer2@gmail.com")
$ws.Range("A9").Style = "Hyperlink"

# --- Column sizing to fit the new content ---
$ws.Columns("A:H").AutoFit()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection, matching the saved cursor position ---
$null = $ws.Range("D12").Select()
